$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Extend the hidden "history" block of date columns from K:AS (11-45) to
# K:AY (11-51). Inserting brand-new columns right at the boundary column
# (AS) makes them inherit AS's column formatting (width 0 / hidden), which
# merges them into the existing hidden range instead of getting a freshly
# computed default width.
# ---------------------------------------------------------------------------
$ws.Columns("AS:AX").Insert() | Out-Null

# The insert pushed the old AS:BE block (13 columns) six places to the
# right, landing on AY:BK. Copy that data back onto AS:BE so every date
# keeps its original column letter, then drop the now-empty columns the
# shift left behind.
$shifted = $ws.Range("AY1:BK18").Value2
$ws.Range("AS1:BE18").Value2 = $shifted
$ws.Range("BF1:BK18").EntireColumn.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Add the new "25-ago" column of data.
# ---------------------------------------------------------------------------
$ws.Range("BF1").Value = "25-ago"
$ws.Range("BF2").Value = 0
$ws.Range("BF3").Value = 12.019981377311693
$ws.Range("BF4").Value = 19.974989025376676
$ws.Range("BF5").Value = 19.470743843545236
$ws.Range("BF6").Value = 0
$ws.Range("BF7").Value = 11.720508748751541
$ws.Range("BF8").Value = 21.73992723887012
$ws.Range("BF9").Value = 9.9739757315545461
$ws.Range("BF10").Value = 1.5685037224286367
$ws.Range("BF11").Value = 11.10824088360085
$ws.Range("BF12").Value = 0
$ws.Range("BF13").Value = 9.8822750163412625
$ws.Range("BF14").Value = 0
$ws.Range("BF15").Value = 0
$ws.Range("BF16").Value = 15.126453384817784
$ws.Range("BF17").Value = 0
$ws.Range("BF18").Value = 0

# ---------------------------------------------------------------------------
# Move the active selection to BH4 (also resets the scrolled-away
# topLeftCell back to the sheet's default top-left corner).
# ---------------------------------------------------------------------------
$ws.Range("BH4").Select() | Out-Null
